$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 85; $row++) {
    $ws.Cells.Item($row, 15).Value = "2022-07-31 20:58:20"
}
